$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Errortype"
$ws.Range("F1").Value = "LOC"

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 18

$ws.Range("F12").Select()
